$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Portefølje_F2026")

# Row 2
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = 508
$ws.Range("E2").Value = 0
$ws.Range("G2").Value = 608
$ws.Range("H2").Value = 50

# Row 3
$ws.Range("E3").Value = 658
$ws.Range("G3").Value = 658
$ws.Range("H3").Value = 0

# Row 4
$ws.Range("D4").Value = 112
$ws.Range("E4").Value = 66
$ws.Range("H4").Value = 480
$ws.Range("J4").Value = 658

# Row 8
$ws.Range("B8").Value = 227
$ws.Range("C8").Value = 376
$ws.Range("E8").Value = 54

# Row 9
$ws.Range("B9").Value = 296
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 389
$ws.Range("I9").ClearContents()

# Row 11
$ws.Range("B11").Value = 369598
$ws.Range("C11").Value = 230852
$ws.Range("D11").Value = 398862
$ws.Range("E11").Value = 299994
